$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.570.05"
$ws.Range("E2").Value = "  -0.56%  "

$ws.Range("D3").Value = "'1.883.38"
$ws.Range("E3").Value = "  -1.38%  "

$ws.Range("D4").Value = "'1.015"
$ws.Range("E4").Value = "  -2.13%  "

$ws.Range("D5").Value = "'317.15"
$ws.Range("E5").Value = "  -1.00%  "

$ws.Range("E6").Value = "  -1.46%  "

$ws.Range("D7").Value = "'0.5129"
$ws.Range("E7").Value = "  -1.68%  "

$ws.Range("D8").Value = "'0.3975"
$ws.Range("E8").Value = "  +0.43%  "

$ws.Range("D9").Value = "'0.08453"
$ws.Range("E9").Value = "  +0.74%  "

$ws.Range("E10").Value = "  -1.80%  "

$ws.Range("D11").Value = "'6.283"
$ws.Range("E11").Value = "  -0.55%  "

$ws.Range("D12").Value = "'1.898.89"
$ws.Range("E12").Value = "  -0.66%  "

$ws.Range("D13").Value = "'20.57"
$ws.Range("E13").Value = "  -0.66%  "

$ws.Range("D14").Value = "'7.298"
$ws.Range("E14").Value = "  -0.41%  "

$ws.Range("D15").Value = "'1.015"
$ws.Range("E15").Value = "  -2.37%  "

$ws.Range("D16").Value = "'0.00001112"
$ws.Range("E16").Value = "  -0.43%  "

$ws.Range("D17").Value = "'91.46"
$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("D18").Value = "'0.06754"
$ws.Range("E18").Value = "  -0.99%  "

$ws.Range("D19").Value = "'17.77"
$ws.Range("E19").Value = "  -1.39%  "

$ws.Range("D20").Value = "'1.014"
$ws.Range("E20").Value = "  -1.35%  "

$ws.Range("D21").Value = "'5.966"
$ws.Range("E21").Value = "  -2.27%  "

$ws.Range("D22").Value = "'28.613.98"
$ws.Range("E22").Value = "  -0.58%  "

$ws.Range("D23").Value = "'11.18"
$ws.Range("E23").Value = "  -0.97%  "

$ws.Range("D24").Value = "'2.272"
$ws.Range("E24").Value = "  -1.00%  "

$ws.Range("D25").Value = "'2.110.95"
$ws.Range("E25").Value = "  -0.78%  "

$ws.Range("D26").Value = "'161.75"
$ws.Range("E26").Value = "  -0.68%  "

$ws.Range("E27").Value = "  -0.82%  "

$ws.Range("D28").Value = "'2.394"
$ws.Range("E28").Value = "  -2.81%  "

$ws.Range("D29").Value = "'127.61"
$ws.Range("E29").Value = "  -0.49%  "

$ws.Range("D30").Value = "'0.1057"
$ws.Range("E30").Value = "  -1.54%  "

$ws.Range("D31").Value = "'1.052"
$ws.Range("E31").Value = "  -0.25%  "

$ws.Range("D32").Value = "'5.796"
$ws.Range("E32").Value = "  -3.69%  "

$ws.Range("D33").Value = "'3.622"
$ws.Range("E33").Value = "  -1.70%  "

$ws.Range("D34").Value = "'0.02448"
$ws.Range("E34").Value = "  -1.05%  "

$ws.Range("E35").Value = "  -1.86%  "

$ws.Range("D36").Value = "'0.2192"
$ws.Range("E36").Value = "  -1.41%  "

$ws.Range("D37").Value = "'8.924"
$ws.Range("E37").Value = "  -5.69%  "

$ws.Range("D38").Value = "'1.269"
$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("D39").Value = "'0.6486"
$ws.Range("E39").Value = "  -1.41%  "

$ws.Range("D40").Value = "'1.193"
$ws.Range("E40").Value = "  -0.49%  "

$ws.Range("D41").Value = "'5.062"
$ws.Range("E41").Value = "  +0.62%  "

$ws.Range("D42").Value = "'11.24"
$ws.Range("E42").Value = "  +0.51%  "

$ws.Range("E43").Value = "  -1.35%  "

$ws.Range("D44").Value = "'0.6096"
$ws.Range("E44").Value = "  -1.47%  "

$ws.Range("D45").Value = "'13.12"
$ws.Range("E45").Value = "  -1.42%  "

$ws.Range("D46").Value = "'3.721"
$ws.Range("E46").Value = "  -1.03%  "

$ws.Range("D47").Value = "'2.024"
$ws.Range("E47").Value = "  +0.07%  "

$ws.Range("D48").Value = "'1.208"
$ws.Range("E48").Value = "  -7.27%  "

$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "'1.216"
$ws.Range("E49").Value = "  -2.21%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'122.89"
$ws.Range("E50").Value = "  +0.14%  "

$ws.Range("D51").Value = "'0.06853"
$ws.Range("E51").Value = "  -2.03%  "

